$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bugfix: fullNormX slider start values (row 7, bern logitX related bugfix).
$ws.Range("T7").Value = "c(-1,.6,.25, 1,-1)"

# Update the funcFormRange (column Q) values for several distribution rows.
$ws.Range("Q5").Value = "c(-6,6)"
$ws.Range("Q6").Value = "c(-6,6)"
$ws.Range("Q7").Value = "c(-8,8)"
$ws.Range("Q8").Value = "c(-2,2)"
$ws.Range("Q9").Value = "c(-2,2)"
$ws.Range("Q10").Value = "c(0,30)"
$ws.Range("Q11").Value = "c(0,30)"
$ws.Range("Q12").Value = "c(0,30)"
$ws.Range("Q13").Value = "c(0,30)"
$ws.Range("Q14").Value = "c(0,30)"
$ws.Range("Q15").Value = "c(0,30)"
$ws.Range("Q16").Value = "c(0,30)"
$ws.Range("Q17").Value = "c(-10,10)"
$ws.Range("Q18").Value = "c(-10,10)"

# Update the selected cell/range shown when the workbook was last saved.
$ws.Range("N12").Select()
